# Generate Report for Handback
#
# Refreshes the handoff/handback timestamps recorded on the per-locale
# report sheets, as if a new handback round had just been generated.
#
# zh-cn sheet: row 2 (9724b7fe-...) gets new "Correspond Handoff Datetime"
#              and "Correspond Handback DateTime" values.
# de-de sheet: row 2 (9724b7fe-...) gets new "Correspond Handoff Datetime"
#              and "Correspond Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-28 08:47:11"
$zhcn.Range("K2").Value = "2016-08-28 08:47:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-28 08:47:15"
$dede.Range("K2").Value = "2016-08-28 08:47:34"
